$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.098.84"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.903.16"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'324.99"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "'0.4609"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").Value = "'0.3880"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "'0.07854"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'0.9891"
$ws.Range("D11").Value = "'21.95"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "1.888.38"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'7.030"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.743"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'0.07025"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'87.97"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'0.000009916"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "29.128.04"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "2.098.07"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'2.088"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "'156.12"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'19.44"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "'5.893"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "'118.56"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "'1.874"
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("D31").Value = "'0.09338"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'0.8939"
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").Value = "'3.136"
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("D36").Value = "'0.05778"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "'1.168"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "'0.02085"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'0.5695"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "'7.637"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("D42").Value = "'0.1803"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'9.702"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Value = "'11.82"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").Value = "'0.5344"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "'0.000002779"
$ws.Range("E46").Value = "  +69.64%  "
$ws.Range("D47").Value = "'2.158"
$ws.Range("E47").Value = "  -6.05%  "
$ws.Range("D48").Value = "'0.06977"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'1.837"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.549"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'112.99"
$ws.Range("E51").Value = "  -0.33%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
